$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 previously held EDINSON VELASQUEZ PALACIN (CC 73578269); replace with MIRELIS VALIENTE MARTINEZ (CC 1002186478)
$ws.Range("C16").Value = "1002186478"
$ws.Range("D16").Value = "MIRELIS VALIENTE MARTINEZ"
$ws.Range("G16").Value = 877803

# Row 17 previously held MIRELIS VALIENTE MARTINEZ (CC 1002186478); replace with EDINSON VELASQUEZ PALACIN (CC 73578269)
$ws.Range("C17").Value = "73578269"
$ws.Range("D17").Value = "EDINSON VELASQUEZ PALACIN"

# Re-fit the columns so the "best fit" widths reflect the new, slightly wider cell content
# (values chosen so the saved column width lands as close as possible to the
# target width on this engine's width model)
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333332
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333336
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
